$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,2).Range.Text = "Max Rice"
$t.Cell(2,2).Range.Text = "Max Rice"
$t.Cell(3,2).Range.Text = "1926372"

$blankAfterTable = $d.Content.Paragraphs.Item(10)
$trailingPara = $d.Content.Paragraphs.Item(12)
$targetRange = $d.Range($blankAfterTable.Range.End, $trailingPara.Range.End)

$combinedXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:r><w:t xml:space="preserve">In this assignment I was supposed to calculate the future population, population change, and determine if there was an increase or decrease in population given the current population, time between death, time between births, and time between immigration, along with how many years into the future we are forecasting for. We accomplished writing a code that will tell us the </w:t></w:r><w:r><w:t>future population, population change, and determine if there was an increase or decrease in population</w:t></w:r><w:r><w:t xml:space="preserve">. We also created an algorithm for our code and an excel sheet to help test our algorithms. To solve this </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>problem</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> we used formulas that could calculate the outputs from the inputs, we also used </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pycharm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and excel applications when working in this project. The key concepts explored was the role mathematical equations and formulas played in creating code</w:t></w:r><w:r><w:t xml:space="preserve"> and how we can use them to create solutions to real world problems. We got the results we were expecting right away with our code but took a while to get what we wanted with excel. Me and my partner both struggled </w:t></w:r><w:r><w:t xml:space="preserve">with operating excel. A challenged I encountered was getting my test cases to work because I did not have excel installed which prevented me from uploading it. I followed the rules by anticipating the hard parts and pushing through them. I overcame them by staying patient and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>trouble shooting</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> my lab partner. I learned how to use equations in excel was something I learned. Working with Jena went </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>well</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> and she helped explain many concepts too me as well as met with me outside class so we could work on this assignment together.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$targetRange.InsertXML($combinedXml)
